# DSE_FH_ScAg_003: Update the requirements according to the project release
#
# Changes the "R" markers to "Q" markers for:
#   - D111:N111  (RQ-SYS.027-9 row)
#   - D121:N125  (RQ-SYS.031 .. RQ-SYS.031-4 rows)
# D111:N111 additionally gets a brand-new bold+underline red Wingdings-2
# style, while D121:N125 reuse the existing "Q" style already used
# elsewhere in the sheet.
# Also updates the active sheet view (frozen pane anchor + selection) to
# reflect where the author was working when the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 111: D111:N111 -> "Q" with a new bold/underline red style ---
# (Font.Color is set before Name/Bold/Underline so the engine never has to
#  materialize a throw-away intermediate font in between.)
$rng111 = $ws.Range("D111:N111")
$rng111.Value = "Q"
$rng111.Font.Color = 255
$rng111.Font.Name = "Wingdings 2"
$rng111.Font.Bold = $true
$rng111.Font.Underline = $true
$rng111.HorizontalAlignment = -4108
$rng111.VerticalAlignment = -4108

# --- Rows 121-125: D:N -> "Q" reusing the existing red Wingdings-2 style ---
$rngBlock = $ws.Range("D121:N125")
$rngBlock.Value = "Q"
$rngBlock.Font.Color = 255
$rngBlock.Font.Name = "Wingdings 2"
$rngBlock.Font.Bold = $true
$rngBlock.Font.Underline = $false
$rngBlock.HorizontalAlignment = -4108
$rngBlock.VerticalAlignment = -4108

# --- Sheet view: move the frozen-pane anchor and the active selection ---
$ws.Application.ActiveWindow.ScrollRow = 59
$ws.Range("D111:N111").Select()
